$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.120.01"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.640.24"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.90"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0635"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.85"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.27"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "1.865.24"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "1.631.05"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.39"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "26.068.41"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.996"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.46"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.59"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.04"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.37"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.67"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.126"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.92"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0500"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.25"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.60"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.909"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").Value = "1.140.03"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.552"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.60"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.68"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.794"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("D45").Value = "1.774.79"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.05"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("D47").Value = "0.0₆0102"
$ws.Range("E47").Value = "  -9.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.48"
$ws.Range("E48").Value = "  +5.39%  "
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.72"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.417"
$ws.Range("E51").Value = "  +0.07%  "
